# UPD checked test cases
# Review comments (column J) added to the Search_Premium test-case sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search_Premium")
$ws.Activate()

# --- Rows 3-11: existing bold column J cells -> turn the font red (reviewer marks) ---
$ws.Range("J3:J11").Font.Color = 255

$ws.Range("J3").Value  = "Лучше в пассиве написать, что context menu is opened…"
$ws.Range("J4").Value  = "Prerequisites должен быть в нашем случае для каждого кейса"
$ws.Range("J5").Value  = "Описать что именно происходит в Expected "
$ws.Range("J7").Value  = "Какой именно dialog window? Что оно спрашивает. Описать"
$ws.Range("J8").Value  = "Проверить что правда сохранилось"
$ws.Range("J9").Value  = "Этот кейс можно обьединить с предыдущим"
$ws.Range("J10").Value = "Prerequisites должен быть в нашем случае для каждого кейса, Лучше в пассиве написать, что context menu is opened…"
$ws.Range("J11").Value = "расписать, нам же надо проверить что документ откроется, а не explorer window"
# J6 keeps the reviewer-red style but no text.

# --- Rows 12-23: brand-new column J cells (non-bold base) -> red reviewer comments ---
$ws.Range("J12").Font.Color = 255
$ws.Range("J13").Font.Color = 255

$ws.Range("J14").Value = "Проверить что правда сохранилось"
$ws.Range("J14").Font.Color = 255
$ws.Range("J14").Font.Bold = $true

$ws.Range("J15").Font.Color = 255

$ws.Range("J16").Value = "Проверить что правда сохранилось"
$ws.Range("J16").Font.Color = 255
$ws.Range("J16").Font.Bold = $true

$ws.Range("J17").Font.Color = 255

$ws.Range("J18").Value = "Проверить что правда сохранилось"
$ws.Range("J18").Font.Color = 255
$ws.Range("J18").Font.Bold = $true

$ws.Range("J19").Font.Color = 255

$ws.Range("J20").Value = "расписать проверку что при изменении параметров они правда меняют что-то на странице"
$ws.Range("J20").Font.Color = 255

$ws.Range("J21").Font.Color = 255
$ws.Range("J22").Font.Color = 255

$ws.Range("J23").Value = "Расписать закрытие пустого файла, с текстом сохраненного, с текстом не сохраненного"
$ws.Range("J23").Font.Color = 255

# --- Misc: restore the cursor position reported by the authored workbook ---
$ws.Range("K10").Select()
